$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-LatticeCell {
    param($cell, $lines)
    $runInner = ""
    for ($i = 0; $i -lt $lines.Length; $i++) {
        if ($i -gt 0) { $runInner += "<w:br/>" }
        $txt = $lines[$i]
        $needsPreserve = ($txt.Length -gt 0) -and (($txt.Substring(0,1) -eq ' ') -or ($txt.Substring($txt.Length-1,1) -eq ' '))
        if ($needsPreserve) {
            $runInner += '<w:t xml:space="preserve">' + $txt + "</w:t>"
        } else {
            $runInner += "<w:t>" + $txt + "</w:t>"
        }
    }
    $runXml = "<w:r><w:rPr><w:sz w:val=" + [char]34 + "32" + [char]34 + "/></w:rPr>" + $runInner + "</w:r>"
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $cell.Range.InsertXML($xmlFrag) | Out-Null
}

Set-LatticeCell $tbl.Cell(1,1) @('46 x 39', '  3    9', '  ----', '4|    |', '6|    |')
Set-LatticeCell $tbl.Cell(1,2) @('80 x 70', '  7    0', '  ----', '8|    |', '0|    |')
Set-LatticeCell $tbl.Cell(1,3) @('76 x 70', '  7    0', '  ----', '7|    |', '6|    |')
Set-LatticeCell $tbl.Cell(2,1) @('37 x 16', '  1    6', '  ----', '3|    |', '7|    |')
Set-LatticeCell $tbl.Cell(2,2) @('14 x 13', '  1    3', '  ----', '1|    |', '4|    |')
Set-LatticeCell $tbl.Cell(2,3) @('84 x 31', '  3    1', '  ----', '8|    |', '4|    |')
Set-LatticeCell $tbl.Cell(3,1) @('75 x 79', '  7    9', '  ----', '7|    |', '5|    |')
Set-LatticeCell $tbl.Cell(3,2) @('26 x 22', '  2    2', '  ----', '2|    |', '6|    |')
Set-LatticeCell $tbl.Cell(3,3) @('47 x 54', '  5    4', '  ----', '4|    |', '7|    |')
Set-LatticeCell $tbl.Cell(4,1) @('54 x 49', '  4    9', '  ----', '5|    |', '4|    |')
Set-LatticeCell $tbl.Cell(4,2) @('85 x 56', '  5    6', '  ----', '8|    |', '5|    |')
Set-LatticeCell $tbl.Cell(4,3) @('45 x 17', '  1    7', '  ----', '4|    |', '5|    |')
Set-LatticeCell $tbl.Cell(5,1) @('49 x 40', '  4    0', '  ----', '4|    |', '9|    |')
Set-LatticeCell $tbl.Cell(5,2) @('22 x 43', '  4    3', '  ----', '2|    |', '2|    |')
Set-LatticeCell $tbl.Cell(5,3) @('84 x 36', '  3    6', '  ----', '8|    |', '4|    |')
